# The regression table gains a new "Intercept" row of coefficients that is
# inserted above the existing "Mining" row (i.e. becomes the new row 2),
# pushing every subsequent row (Mining, Partisanship, Deregulated, R-squared,
# N, State FE) down by one. The sheet's dimension grows from A1:G7 to A1:G8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 2; Excel shifts rows 2-7 down to 3-8
# and updates the sheet dimension automatically.
$ws.Rows.Item(2).Insert()

# The inserted row inherits formatting from the row above it (the bold
# header row), so strip that back to the workbook's default formatting
# before applying the correct look to the new row-label cell (A2).
$ws.Range("A2:G2").ClearFormats()

# Re-apply the same look used by the other row-label cells in column A
# (bold, centered/top aligned, thin box border) to the new "Intercept" label.
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4160
$ws.Range("A2").Borders.LineStyle = 1
$ws.Range("A2").Borders.Weight = 2

# Fill in the new row's contents. Columns C, E and G stay blank for this
# row (those model specifications have no separately estimated intercept),
# matching the existing pattern used by the "Deregulated" row above.
$ws.Range("A2").Value = "Intercept"
$ws.Range("B2").Value = "0.244**`n (0.096)"
$ws.Range("D2").Value = "0.343***`n (0.068)"
$ws.Range("F2").Value = "0.511***`n (0.052)"

# Setting multi-line text makes Excel apply a custom "autofit" row height;
# explicitly auto-fitting afterwards restores the standard row height so
# the row doesn't carry an explicit custom height.
$ws.Rows.Item(2).EntireRow.AutoFit() | Out-Null
